$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the style of the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the time_taken values for each data row
$times = @(
    "2021-10-05 13:41:38.993660",
    "2021-10-05 13:41:38.993671",
    "2021-10-05 13:41:38.993675",
    "2021-10-05 13:41:38.993678",
    "2021-10-05 13:41:38.993682",
    "2021-10-05 13:41:38.993685"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
